# Fruta / hortaliza, semanal
# Insert two new daily price records for "Agrícola del Norte S.A. de Arica - Uva"
# right after the existing row 54, pushing the following rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 55 and 56 (old row 55 data shifts to row 57, etc.)
$ws.Rows("55:56").Insert()

# --- New row 55: Flame Seedless / Primera ---
$ws.Range("A55").Value2 = 1
$ws.Range("B55").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C55").Value2 = 'Arica y Parinacota'
$ws.Range("D55").Value2 = 44588
$ws.Range("E55").Value2 = 15
$ws.Range("F55").Value2 = 'Fruta'
$ws.Range("G55").Value2 = 100109
$ws.Range("H55").Value2 = 'Uva'
$ws.Range("I55").Value2 = 100109001
$ws.Range("J55").Value2 = 'Uva'
$ws.Range("K55").Value2 = 'Flame Seedless'
$ws.Range("L55").Value2 = 'Primera'
$ws.Range("M55").Value2 = 270
$ws.Range("N55").Value2 = 21000
$ws.Range("O55").Value2 = 22000
$ws.Range("P55").Value2 = 21500
$ws.Range("Q55").Value2 = '$/caja 25 kilos'
$ws.Range("R55").Value2 = 'Región de Coquimbo'
$ws.Range("S55").Value2 = 860
$ws.Range("T55").Value2 = 25

# --- New row 56: Superior Seedless / Segunda ---
$ws.Range("A56").Value2 = 1
$ws.Range("B56").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C56").Value2 = 'Arica y Parinacota'
$ws.Range("D56").Value2 = 44588
$ws.Range("E56").Value2 = 15
$ws.Range("F56").Value2 = 'Fruta'
$ws.Range("G56").Value2 = 100109
$ws.Range("H56").Value2 = 'Uva'
$ws.Range("I56").Value2 = 100109001
$ws.Range("J56").Value2 = 'Uva'
$ws.Range("K56").Value2 = 'Superior Seedless'
$ws.Range("L56").Value2 = 'Segunda'
$ws.Range("M56").Value2 = 250
$ws.Range("N56").Value2 = 17000
$ws.Range("O56").Value2 = 18000
$ws.Range("P56").Value2 = 17500
$ws.Range("Q56").Value2 = '$/bandeja 18 kilos'
$ws.Range("R56").Value2 = 'Región de Coquimbo'
$ws.Range("S56").Value2 = 972
$ws.Range("T56").Value2 = 18
